# Add the "2022-Q1" holdings sheet (before "总计") and update the "总计"
# summary sheet with a new leading row for 2022-Q1.
#
# NOTE: sheet object references returned by Worksheets.Item(...) track a
# *position*, not a stable identity - once Worksheets.Add() shifts indices
# around, any previously-grabbed reference can silently start pointing at a
# different sheet. So every sheet we touch is re-fetched by name right
# before it is used.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q1" sheet right before the "总计" sheet.
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($wb.Worksheets.Item("总计"))
$newSheet.Name = "2022-Q1"

# A same-layout sheet ("2021-Q4") supplies the header/index cell style (s=2:
# bold, centered, thin border) and a plain (unstyled) cell, so the new
# sheet matches the others exactly.
$styleSrc = $wb.Worksheets.Item("2021-Q4")

# Header row (B1:H1) - bold/centered/bordered style.
$styleSrc.Range("B1:H1").Copy()
$wb.Worksheets.Item("2022-Q1").Range("B1:H1").PasteSpecial(-4122)
$wb.Worksheets.Item("2022-Q1").Range("B1").Value = "基金代码"
$wb.Worksheets.Item("2022-Q1").Range("C1").Value = "基金名称"
$wb.Worksheets.Item("2022-Q1").Range("D1").Value = "基金规模"
$wb.Worksheets.Item("2022-Q1").Range("E1").Value = "股票总仓位"
$wb.Worksheets.Item("2022-Q1").Range("F1").Value = "仓位占比"
$wb.Worksheets.Item("2022-Q1").Range("G1").Value = "持有市值(亿元)"
$wb.Worksheets.Item("2022-Q1").Range("H1").Value = "仓位排名"

# Fund holdings data rows.
$fundRows = @(
    @("166002", "中欧新蓝筹混合 -A", "129.77", "77.81", "2.21", "2.8679", 10),
    @("001117", "中欧精选灵活配置定期开放混合A", "59.37", "86.95", "2.36", "1.4011", 10),
    @("001890", "中欧精选灵活配置定期开放混合E", "59.37", "86.95", "2.36", "1.4011", 10),
    @("166023", "中欧瑞丰灵活配置混合（LOF）A", "32.40", "85.04", "3.27", "1.0595", 8),
    @("000314", "招商瑞丰灵活配置混合A", "14.01", "40.96", "1.80", "0.2522", 8),
    @("001427", "招商丰泽灵活配置混合A", "13.94", "41.14", "1.64", "0.2286", 9),
    @("004237", "中欧新蓝筹混合 -C", "9.82", "77.81", "2.21", "0.2170", 10),
    @("002389", "招商安德灵活配置混合A", "8.03", "44.29", "1.89", "0.1518", 7),
    @("002819", "招商丰美灵活配置混合A", "7.83", "29.10", "1.46", "0.1143", 7),
    @("002017", "招商瑞丰灵活配置混合C", "4.66", "40.96", "1.80", "0.0839", 8),
    @("013393", "信达澳银价值精选混合A", "3.61", "81.31", "1.84", "0.0664", 10),
    @("002390", "招商安德灵活配置混合C", "2.55", "44.29", "1.89", "0.0482", 7),
    @("004740", "中欧瑞丰灵活配置混合（LOF）C", "1.28", "85.04", "3.27", "0.0419", 8),
    @("001446", "招商丰泽灵活配置混合C", "1.99", "41.14", "1.64", "0.0326", 9),
    @("001885", "中欧新蓝筹混合 -E", "1.41", "77.81", "2.21", "0.0312", 10),
    @("002820", "招商丰美灵活配置混合C", "0.73", "29.10", "1.46", "0.0107", 7),
    @("013394", "信达澳银价值精选混合C", "0.37", "81.31", "1.84", "0.0068", 10)
)

$r = 2
foreach ($row in $fundRows) {
    $fundSheet = $wb.Worksheets.Item("2022-Q1")

    # Column A: numeric row index (0-based), same bold/bordered style as header.
    $styleSrc.Range("A2").Copy()
    $fundSheet.Cells.Item($r, 1).PasteSpecial(-4122)
    $fundSheet.Cells.Item($r, 1).Value = ($r - 2)

    # Columns B-G: always stored as text (fund codes keep leading zeros,
    # decimal figures keep their exact printed form).
    for ($c = 2; $c -le 7; $c++) {
        $cell = $fundSheet.Cells.Item($r, $c)
        $cell.NumberFormat = "@"
        $cell.Value = [string]$row[$c - 2]
        $styleSrc.Range("C2").Copy()
        $cell.PasteSpecial(-4122)
    }

    # Column H: numeric rank.
    $fundSheet.Cells.Item($r, 8).Value = $row[6]

    $r = $r + 1
}

$wb.Worksheets.Item("2022-Q1").Range("A1").Select()

# ---------------------------------------------------------------------
# 2. Update the "总计" summary sheet: add a new leading row for 2022-Q1
#    and shift the existing rows down, re-numbering the index column.
# ---------------------------------------------------------------------
$summaryRows = @(
    @("2022-Q1", 17, 8.02),
    @("2021-Q4", 1, 0.05),
    @("2021-Q3", 18, 12.51),
    @("2021-Q1", 2, 0.01)
)

$r = 2
foreach ($row in $summaryRows) {
    $totalSheet = $wb.Worksheets.Item("总计")
    $totalSheet.Range("A2").Copy()
    $totalSheet.Cells.Item($r, 1).PasteSpecial(-4122)
    $totalSheet.Cells.Item($r, 1).Value = ($r - 2)
    $totalSheet.Cells.Item($r, 2).Value = $row[0]
    $totalSheet.Cells.Item($r, 3).Value = $row[1]
    $totalSheet.Cells.Item($r, 4).Value = $row[2]
    $r = $r + 1
}

$wb.Worksheets.Item("总计").Range("A1").Select()
